$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 12346682
$ws.Range("I6").Value = 41666856
$ws.Range("J6").Value = 1346.579
$ws.Range("K6").Value = 125000568
$ws.Range("L6").Value = 4039.737
$ws.Range("M6").Value = -125000456
$ws.Range("N6").Value = -4263.737
$ws.Range("H8").Value = 15151745
$ws.Range("I8").Value = 25641208
$ws.Range("K8").Value = 76923624
$ws.Range("M8").Value = -76923485
$ws.Range("H51").Value = 40719.465
$ws.Range("I51").Value = 7790.5
$ws.Range("K51").Value = 7790.5
$ws.Range("M51").Value = -7306.5
$ws.Range("H55").Value = 94.25
$ws.Range("J55").Value = 53
$ws.Range("L55").Value = 53
$ws.Range("N55").Value = -481
$ws.Range("H98").Value = 27210.092
$ws.Range("I98").Value = 32013.445
$ws.Range("J98").Value = 5595
$ws.Range("K98").Value = 32013.445
$ws.Range("L98").Value = 5595
$ws.Range("M98").Value = -30515.445
$ws.Range("N98").Value = -8591
$ws.Range("H111").Value = 920.5333000000001
$ws.Range("I111").Value = 720.7273
$ws.Range("K111").Value = 2162.1819
$ws.Range("M111").Value = 904.8181
$ws.Range("H112").Value = 2279.9375
$ws.Range("I112").Value = 4000
$ws.Range("J112").Value = 1883
$ws.Range("K112").Value = 12000
$ws.Range("L112").Value = 5649
$ws.Range("M112").Value = -10892
$ws.Range("N112").Value = -7865
$ws.Range("H122").Value = 27210.092
$ws.Range("I122").Value = 32013.445
$ws.Range("J122").Value = 5595
$ws.Range("K122").Value = 96040.33499999999
$ws.Range("L122").Value = 16785
$ws.Range("M122").Value = -93590.33499999999
$ws.Range("N122").Value = -21685
$ws.Range("H125").Value = 11750
$ws.Range("J125").Value = 5000
$ws.Range("L125").Value = 45000
$ws.Range("N125").Value = -49920
$ws.Range("H132").Value = 1821819.2
$ws.Range("I132").Value = 3849.25
$ws.Range("K132").Value = 11547.75
$ws.Range("M132").Value = -9017.75

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10907.76
$ws.Range("I32").Value = 10306.841
$ws.Range("K32").Value = 10306.841
$ws.Range("M32").Value = -10019.841
$ws.Range("H45").Value = 189441.81
$ws.Range("I45").Value = 408858.2
$ws.Range("K45").Value = 408858.2
$ws.Range("M45").Value = -408481.2
$ws.Range("H61").Value = 7586.4883
$ws.Range("I61").Value = 8089.387
$ws.Range("K61").Value = 8089.387
$ws.Range("M61").Value = -7877.387
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").ClearContents()
$ws.Range("H74").Value = 1635.3334
$ws.Range("I74").Value = 604.6875
$ws.Range("K74").Value = 604.6875
$ws.Range("M74").Value = 269.3125
$ws.Range("H77").Value = 1635.3334
$ws.Range("I77").Value = 604.6875
$ws.Range("K77").Value = 3023.4375
$ws.Range("M77").Value = 1344.5625
$ws.Range("H132").Value = 2015.037
$ws.Range("I132").Value = 1246.1
$ws.Range("J132").Value = 4212
$ws.Range("K132").Value = 3738.3
$ws.Range("L132").Value = 12636
$ws.Range("M132").Value = -1208.3
$ws.Range("N132").Value = -17696
$ws.Range("H136").Value = 7586.4883
$ws.Range("I136").Value = 8089.387
$ws.Range("K136").Value = 24268.161
$ws.Range("M136").Value = -21718.161

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 40000
$ws.Range("J62").Value = 40000
$ws.Range("L62").Value = 40000
$ws.Range("N62").Value = -41372
$ws.Range("H65").Value = 40000
$ws.Range("J65").Value = 40000
$ws.Range("L65").Value = 120000
$ws.Range("N65").Value = -126864
$ws.Range("H96").Value = 16556.916
$ws.Range("I96").Value = 11243.909
$ws.Range("K96").Value = 11243.909
$ws.Range("M96").Value = -8497.909
$ws.Range("H134").Value = 2082.2
$ws.Range("I134").Value = 1288.3334
$ws.Range("K134").Value = 3865.0002
$ws.Range("M134").Value = -1330.0002

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4733.34
$ws.Range("I31").Value = 4854.2446
$ws.Range("J31").Value = 4053.25
$ws.Range("K31").Value = 4854.2446
$ws.Range("L31").Value = 4053.25
$ws.Range("M31").Value = -4559.2446
$ws.Range("N31").Value = -4643.25
$ws.Range("H34").Value = 4733.34
$ws.Range("I34").Value = 4854.2446
$ws.Range("J34").Value = 4053.25
$ws.Range("K34").Value = 4854.2446
$ws.Range("L34").Value = 4053.25
$ws.Range("M34").Value = -4652.2446
$ws.Range("N34").Value = -4457.25
$ws.Range("H58").Value = 2650.9033
$ws.Range("I58").Value = 1883.4762
$ws.Range("K58").Value = 1883.4762
$ws.Range("M58").Value = -1680.4762
$ws.Range("H99").Value = 3420206
$ws.Range("I99").Value = 5810990
$ws.Range("K99").Value = 5810990
$ws.Range("M99").Value = -5809492
$ws.Range("H107").Value = 7732.5884
$ws.Range("I107").Value = 10346.25
$ws.Range("K107").Value = 10346.25
$ws.Range("M107").Value = -8426.25
$ws.Range("H126").Value = 3420206
$ws.Range("I126").Value = 5810990
$ws.Range("K126").Value = 17432970
$ws.Range("M126").Value = -17430500
$ws.Range("H132").Value = 1888.6061
$ws.Range("I132").Value = 1687.9032
$ws.Range("K132").Value = 5063.7096
$ws.Range("M132").Value = -2533.7096
$ws.Range("H136").Value = 2650.9033
$ws.Range("I136").Value = 1883.4762
$ws.Range("K136").Value = 5650.4286
$ws.Range("M136").Value = -3100.4286
$ws.Range("H141").Value = 178679.08
$ws.Range("J141").Value = 188768.95
$ws.Range("L141").Value = 188768.95
$ws.Range("N141").Value = -199128.95

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 3994.4893
$ws.Range("J122").Value = 5070.0884
$ws.Range("L122").Value = 45630.7956
$ws.Range("N122").Value = -50530.7956

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 40061
$ws.Range("J57").Value = 40061
$ws.Range("L57").Value = 40061
$ws.Range("N57").Value = -41701
$ws.Range("H122").Value = 9082.392
$ws.Range("I122").Value = 6458.8237
$ws.Range("J122").Value = 16515.834
$ws.Range("K122").Value = 19376.4711
$ws.Range("L122").Value = 49547.50199999999
$ws.Range("M122").Value = -16926.4711
$ws.Range("N122").Value = -54447.50199999999
$ws.Range("H132").Value = 2036.9487
$ws.Range("I132").Value = 1898.4857
$ws.Range("J132").Value = 3248.5
$ws.Range("K132").Value = 5695.4571
$ws.Range("L132").Value = 9745.5
$ws.Range("M132").Value = -3165.4571
$ws.Range("N132").Value = -14805.5

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2485125.5
$ws.Range("I132").Value = 2980550.5
$ws.Range("K132").Value = 8941651.5
$ws.Range("M132").Value = -8939121.5

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 275030000
$ws.Range("J92").Value = 275030000
$ws.Range("L92").Value = 275030000
$ws.Range("N92").Value = -275034992
$ws.Range("H132").Value = 9714.888999999999
$ws.Range("I132").Value = 10378.211
$ws.Range("J132").Value = 6114
$ws.Range("K132").Value = 31134.633
$ws.Range("L132").Value = 18342
$ws.Range("M132").Value = -28604.633
$ws.Range("N132").Value = -23402
$ws.Range("H136").Value = 286535.03
$ws.Range("I136").Value = 291828.16
$ws.Range("K136").Value = 875484.48
$ws.Range("M136").Value = -872934.48
$ws.Range("H139").Value = 107367.75
$ws.Range("J139").Value = 73157
$ws.Range("L139").Value = 73157
$ws.Range("N139").Value = -83437
